# Add "Toisen asteen pohjakoulutus suoritettu" column.
#
# A brand-new column is inserted right before the existing "Pohjakoulutus
# maa (toinen aste)" column (which was column W and becomes column X).
# The new column gets a header in row 1 and a "Kyllä" value in row 2 (the
# single data row), matching the style ids already used by the
# neighbouring header/data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at W; everything from W onward (formerly
# "Pohjakoulutus maa (toinen aste)" in W) shifts right to X.
$ws.Columns("W").Insert() | Out-Null

# Fill in the new column's header and the single data row's value.
$ws.Range("W1").Value = "Toisen asteen pohjakoulutus suoritettu"
$ws.Range("W2").Value = "Kyllä"

# The author's selection ended up on the newly added data cell.
$ws.Range("W2").Select() | Out-Null
